$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 values changed: G13->H3 labels shifted (WW52.2 -> WW1.1, 5PM-5PM -> 12PM-12PM)
# Update H3 before F3 so the shared-string table order matches the target workbook.
$ws.Range("H3").Value = "WW1.1"
$ws.Range("F3").Value = "12PM-12PM"

# Update the active cell/selection on the sheet to match the saved view state.
$ws.Range("J8").Select()
